$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.700.09'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.877.67'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('ZZ1').Formula = '="1.0000"'
$ws.Range('ZZ1').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('ZZ1').Formula = '="237.36"'
$ws.Range('ZZ1').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('ZZ1').Formula = '="0.9999"'
$ws.Range('ZZ1').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('ZZ1').Formula = '="0.4733"'
$ws.Range('ZZ1').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('ZZ1').Formula = '="0.2819"'
$ws.Range('ZZ1').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +3.17%  '
$ws.Range('ZZ1').Formula = '="0.06499"'
$ws.Range('ZZ1').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +3.44%  '
$ws.Range('ZZ1').Formula = '="18.57"'
$ws.Range('ZZ1').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +14.11%  '
$ws.Range('D11').Value = '1.880.79'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('ZZ1').Formula = '="0.07572"'
$ws.Range('ZZ1').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('ZZ1').Formula = '="95.29"'
$ws.Range('ZZ1').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +13.72%  '
$ws.Range('ZZ1').Formula = '="5.071"'
$ws.Range('ZZ1').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('ZZ1').Formula = '="0.6487"'
$ws.Range('ZZ1').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +4.40%  '
$ws.Range('ZZ1').Formula = '="303.60"'
$ws.Range('ZZ1').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +32.91%  '
$ws.Range('D17').Value = '30.687.65'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('ZZ1').Formula = '="13.06"'
$ws.Range('ZZ1').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +5.79%  '
$ws.Range('ZZ1').Formula = '="0.9990"'
$ws.Range('ZZ1').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('ZZ1').Formula = '="0.000007530"'
$ws.Range('ZZ1').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').Value = '2.124.85'
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('ZZ1').Formula = '="0.9998"'
$ws.Range('ZZ1').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('ZZ1').Formula = '="5.136"'
$ws.Range('ZZ1').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('ZZ1').Formula = '="6.138"'
$ws.Range('ZZ1').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +4.69%  '
$ws.Range('ZZ1').Formula = '="168.92"'
$ws.Range('ZZ1').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('ZZ1').Formula = '="9.214"'
$ws.Range('ZZ1').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('ZZ1').Formula = '="19.65"'
$ws.Range('ZZ1').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +10.43%  '
$ws.Range('ZZ1').Formula = '="1.944"'
$ws.Range('ZZ1').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +3.91%  '
$ws.Range('ZZ1').Formula = '="0.1055"'
$ws.Range('ZZ1').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('ZZ1').Formula = '="4.154"'
$ws.Range('ZZ1').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('E32').Value = '  +3.36%  '
$ws.Range('ZZ1').Formula = '="0.05042"'
$ws.Range('ZZ1').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('ZZ1').Formula = '="1.169"'
$ws.Range('ZZ1').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('ZZ1').Formula = '="0.7181"'
$ws.Range('ZZ1').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +0.95%  '
$ws.Range('ZZ1').Formula = '="2.707"'
$ws.Range('ZZ1').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('ZZ1').Formula = '="0.01912"'
$ws.Range('ZZ1').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('ZZ1').Formula = '="2.705"'
$ws.Range('ZZ1').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('ZZ1').Formula = '="2.043"'
$ws.Range('ZZ1').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +6.01%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('ZZ1').Formula = '="0.8959"'
$ws.Range('ZZ1').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('ZZ1').Formula = '="107.02"'
$ws.Range('ZZ1').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('ZZ1').Formula = '="0.9997"'
$ws.Range('ZZ1').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('ZZ1').Formula = '="0.4185"'
$ws.Range('ZZ1').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +4.16%  '
$ws.Range('ZZ1').Formula = '="5.577"'
$ws.Range('ZZ1').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('ZZ1').Formula = '="7.303"'
$ws.Range('ZZ1').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('ZZ1').Formula = '="64.86"'
$ws.Range('ZZ1').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +6.63%  '
$ws.Range('ZZ1').Formula = '="8.939"'
$ws.Range('ZZ1').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +4.25%  '
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('ZZ1').Formula = '="34.54"'
$ws.Range('ZZ1').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +4.06%  '
$ws.Range('ZZ1').Formula = '="0.05585"'
$ws.Range('ZZ1').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('ZZ1').Formula = '="1.375"'
$ws.Range('ZZ1').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +1.78%  '
$ws.Range('ZZ1').ClearContents()
